# Add a new "test_suite" worksheet as the first (left-most) sheet in the
# workbook, containing the Test-Case-ID / Run-mode table, as described by
# the commit "Setting up run modes for Test Suites (Skipping Open Account Test)".

$wb = $excel.ActiveWorkbook

# Create the new worksheet and move it to the very first position in the
# workbook (Worksheets.Add defaults to inserting before the active sheet,
# but we explicitly pass "Before" to be safe).
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($firstSheet)
$newSheet.Name = "test_suite"

# Populate the header row.
$newSheet.Range("A1").Value = "TCID"
$newSheet.Range("B1").Value = "Runmode"

# Populate the test case rows - all tests run ("Y") except OpenAccountTest,
# which is skipped ("N").
$newSheet.Range("A2").Value = "BankManagerLoginTest"
$newSheet.Range("B2").Value = "Y"

$newSheet.Range("A3").Value = "AddCustomerTest"
$newSheet.Range("B3").Value = "Y"

$newSheet.Range("A4").Value = "OpenAccountTest"
$newSheet.Range("B4").Value = "N"

# Make this new first sheet the active / selected tab, and set the
# selection to B4 to match the authored workbook state.
$newSheet.Activate()
$newSheet.Range("B4").Select()
